# Lecture partielle de l'EDT M1 MIAGE.
# Shift the displayed schedule dates forward (from the 2023 week set to the
# 2026 week set) and update the corresponding day-of-week labels to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: jeudi 04/05/2023 -> lundi 04/05/2026
$ws.Range("A2").Value = 46146.0
$ws.Range("B2").Value = "lundi"

# Row 4: dimanche 07/05/2023 -> jeudi 07/05/2026
$ws.Range("A4").Value = 46149.0
$ws.Range("B4").Value = "jeudi"

# Row 7: jeudi 11/05/2023 -> lundi 11/05/2026
$ws.Range("A7").Value = 46153.0
$ws.Range("B7").Value = "lundi"

# Row 10: samedi 20/05/2023 -> mercredi 20/05/2026
$ws.Range("A10").Value = 46162.0
$ws.Range("B10").Value = "mercredi"
